$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.965.77"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "'3.310.08"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'186.84"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'557.97"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").Value = "'3.302.35"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "'0.584"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "'47.66"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "'0.0000270"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'629.03"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("D16").Value = "'3.835.28"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "'18.12"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").Value = "'65.946.12"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("D20").Value = "'3.309.88"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = "  -4.28%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'18.26"
$ws.Range("E23").Value = "  +6.51%  "
$ws.Range("D24").Value = "'101.79"
$ws.Range("E24").Value = "  +7.12%  "
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "'6.04"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").Value = "'9.54"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "'8.65"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "'30.24"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").Value = "'4.00"
$ws.Range("E32").Value = "  +4.17%  "
$ws.Range("D33").Value = "'6.36"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "'11.06"
$ws.Range("D35").Value = "'552.97"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "'3.851.58"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("D38").Value = "'57.55"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'0.0₃0732"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "'33.80"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("D42").Value = "'3.28"
$ws.Range("E42").Value = "  -3.96%  "
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").Value = "'3.24"
$ws.Range("E45").Value = "  -14.61%  "
$ws.Range("E46").Value = "  -5.64%  "
$ws.Range("D47").Value = "'0.0418"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "'3.27"
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'2.60"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.04%  "
